$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42654.743668981479
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"

$ws.Range("D3").Value = 75.5

$ws.Range("G3").Value = $true
$ws.Range("G3").NumberFormat = "m/d/yy h:mm"
